$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new column AA (year 2023) that mirrors the existing
# column Z (year 2022) formatting for rows 4-16. Copy formats first (so we
# reuse/derive styles the same way Excel would), then write the new values.
for ($r = 4; $r -le 16; $r++) {
    $ws.Range("Z$r").Copy()
    $ws.Range("AA$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("AA4").Value = 2023
$ws.Range("AA5").Value = 44.2
$ws.Range("AA6").Value = 50.4
$ws.Range("AA7").Value = 40.6
$ws.Range("AA8").Value = 57.2
$ws.Range("AA9").Value = 31
$ws.Range("AA10").Value = 49.7
$ws.Range("AA11").Value = 51
$ws.Range("AA12").Value = 29.4
$ws.Range("AA13").Value = 29.9
$ws.Range("AA14").Value = 56.3
$ws.Range("AA15").Value = 62.5
$ws.Range("AA16").Value = 34.9

# The view used to be scrolled so column B was the left-most visible column
# with AA4 selected; restore the default (top-left = A1, nothing special
# selected) now that AA is a normal, visible part of the used range.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A1").Select() | Out-Null
